$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 15716.53094379244
$ws.Range("D3").Value = 781.2597083856616

$ws.Range("B4").Value = 6690.590462644857
$ws.Range("D4").Value = 530.3026332473772

$ws.Range("B5").Value = 2544.015

$ws.Range("B6").Value = 10516.012
$ws.Range("D6").Value = 202.003

$ws.Range("B7").Value = 14086.03150000001
$ws.Range("D7").Value = 1040

$ws.Range("B8").Value = 21537.028
$ws.Range("D8").Value = 1080

$ws.Range("B9").Value = 33915.31300000004
$ws.Range("D9").Value = 6088.003

$ws.Range("F10").Value = 8260608664.057012

$ws.Range("G11").Value = 0.8115588775488013

$ws.Range("F12").Value = 410630102.7280001
$ws.Range("G12").Value = 0.04970942450218049

$ws.Range("G13").Value = 0.1387316979490182
